$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '202111 45784 Holy Rosary Primary SchoolWhite Hills Outbreak'
$ws.Cells.Item(2, 2).Value = 27
$ws.Cells.Item(3, 1).Value = '3321 Rochester and Elmore District Health Service Yalunkan Aged Care Hostel RochesterOutbreak'
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(4, 1).Value = '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North Outbreak'
$ws.Cells.Item(4, 2).Value = 12
$ws.Cells.Item(5, 1).Value = '3600 Belvedere Aged Care Noble Park Outbreak'
$ws.Cells.Item(5, 2).Value = 20
$ws.Cells.Item(6, 1).Value = '3601 Baptcare Westhaven community Outbreak'
$ws.Cells.Item(6, 2).Value = 26
$ws.Cells.Item(7, 1).Value = '3653 Fronditha Thalpori St Albans Aged Care Outbreak'
$ws.Cells.Item(7, 2).Value = 33
$ws.Cells.Item(8, 1).Value = '44098 Stawell Primary School Outbreak'
$ws.Cells.Item(8, 2).Value = 14
$ws.Cells.Item(9, 1).Value = '44121 Wallan Primary School Wallan Outbreak'
$ws.Cells.Item(9, 2).Value = 21
$ws.Cells.Item(10, 1).Value = '44165 Greenvale Primary School Outbreak'
$ws.Cells.Item(10, 2).Value = 34
$ws.Cells.Item(11, 1).Value = '44234 Lucknow Primary School Bairnsdale Outbreak'
$ws.Cells.Item(11, 2).Value = 42
$ws.Cells.Item(12, 1).Value = '44495 Lakes Entrance Primary School Outbreak'
$ws.Cells.Item(12, 2).Value = 10
$ws.Cells.Item(13, 1).Value = '44667 Beaumaris Primary School Beaumaris Outbreak'
$ws.Cells.Item(13, 2).Value = 22
$ws.Cells.Item(14, 1).Value = '44811 Dandenong North Primary SchoolDandenong Outbreak'
$ws.Cells.Item(14, 2).Value = 27
$ws.Cells.Item(15, 1).Value = '44853 St Albans North Primary School Outbreak'
$ws.Cells.Item(15, 2).Value = 11
$ws.Cells.Item(16, 1).Value = '44865 Parktone Primary School Parkdale Outbreak'
$ws.Cells.Item(16, 2).Value = 13
$ws.Cells.Item(17, 1).Value = '44891 Cranbourne Park Primary School Cranbourne Outbreak'
$ws.Cells.Item(17, 2).Value = 20
$ws.Cells.Item(18, 1).Value = '45158 Rowellyn Park Primary School Carrum Downs Outbreak'
$ws.Cells.Item(18, 2).Value = 13
$ws.Cells.Item(19, 1).Value = '45248 Brookside P-9 College Caroline Springs Outbreak'
$ws.Cells.Item(19, 2).Value = 14
$ws.Cells.Item(20, 1).Value = '45249 Creekside K-9 College Caroline SpringsOutbreak'
$ws.Cells.Item(20, 2).Value = 16
$ws.Cells.Item(21, 1).Value = '45569 Nhill College Nhill Outbreak'
$ws.Cells.Item(21, 2).Value = 15
$ws.Cells.Item(22, 1).Value = '4574 Village Glen Aged Care Residences Mornington Outbreak'
$ws.Cells.Item(22, 2).Value = 17
$ws.Cells.Item(23, 1).Value = '45836 St Joseph''s Primary School Sorrento Outbreak'
$ws.Cells.Item(23, 2).Value = 16
$ws.Cells.Item(24, 1).Value = '45967 St Clement of Rome School Bulleen Outbreak'
$ws.Cells.Item(24, 2).Value = 10
$ws.Cells.Item(25, 1).Value = '46037 Nazareth Catholic Primary SchoolGrovedale Outbreak'
$ws.Cells.Item(25, 2).Value = 29
$ws.Cells.Item(26, 1).Value = '46050 Our Lady''s Catholic Primary School Craigieburn Outbreak'
$ws.Cells.Item(26, 2).Value = 30
$ws.Cells.Item(27, 1).Value = '46125 Our Lady of the Southern Cross Primary School Manor Lakes Outbreak'
$ws.Cells.Item(27, 2).Value = 29
$ws.Cells.Item(28, 1).Value = '46190 Haileybury Brighton East Outbreak'
$ws.Cells.Item(28, 2).Value = 13
$ws.Cells.Item(29, 1).Value = '46215 Yeshivah primary College St Kilda East Outbreak'
$ws.Cells.Item(29, 2).Value = 12
$ws.Cells.Item(30, 1).Value = '46276 Hillcrest Christian College Clyde NorthOutbreak'
$ws.Cells.Item(30, 2).Value = 18
$ws.Cells.Item(31, 1).Value = '46328 Ilim College Dallas Main Campus Dallas Oct Outbreak'
$ws.Cells.Item(31, 2).Value = 30
$ws.Cells.Item(32, 1).Value = '46376 Yesodei HaTorah College Elwood Outbreak'
$ws.Cells.Item(32, 2).Value = 11
$ws.Cells.Item(33, 1).Value = '46390 Al Siraat College Epping Outbreak'
$ws.Cells.Item(33, 2).Value = 32
$ws.Cells.Item(34, 1).Value = '50395 St Francis of Assisi Catholic PrimarySchool Tarneit Outbreak'
$ws.Cells.Item(34, 2).Value = 11
$ws.Cells.Item(35, 1).Value = '50681 Broadmeadows Special Developmental School Broadmeadows Outbreak'
$ws.Cells.Item(35, 2).Value = 12
$ws.Cells.Item(36, 1).Value = '52380 Al Iman College Melton South Outbreak'
$ws.Cells.Item(36, 2).Value = 14
$ws.Cells.Item(37, 1).Value = '52473 John Henry Primary School PakenhamOutbreak'
$ws.Cells.Item(37, 2).Value = 19
$ws.Cells.Item(38, 1).Value = 'Adass Israel School Elsternwick Outbreak'
$ws.Cells.Item(38, 2).Value = 10
$ws.Cells.Item(39, 1).Value = 'Hamilton Country Music Festival Hamilton Golf Club Hamilton Outbreak'
$ws.Cells.Item(39, 2).Value = 27
$ws.Cells.Item(40, 1).Value = 'Melton Willows Melton Outbreak'
$ws.Cells.Item(40, 2).Value = 11
$ws.Cells.Item(41, 1).Value = 'St Brendans Primary School Shepparton Outbreak'
$ws.Cells.Item(41, 2).Value = 12
